# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columna I (municipio-nombre): pasa de ser medida (iaest-measure:municipio-nombre,
#   xsd:int) a ser dimension (sdmx-dimension:refArea, URI-Municipio).
# Columna K (sexo): pasa de ser dimension (iaest-dimension:sexo, skos:Concept,
#   mapping-sexo.xlsx) a ser medida (iaest-measure:sexo, xsd:int, sin mapping).
# Columna L (direccion-provincial-nombre): pasa de ser dimension
#   (sdmx-dimension:refArea, URI-Provincia) a ser medida
#   (iaest-measure:direccion-provincial-nombre, xsd:int).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columna I - municipio-nombre
$ws.Range("I2").Value = "sdmx-dimension:refArea"
$ws.Range("I3").Value = "dim"
$ws.Range("I4").Value = "URI-Municipio"

# Columna K - sexo
$ws.Range("K2").Value = "iaest-measure:sexo"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("K5").Clear()

# Columna L - direccion-provincial-nombre
$ws.Range("L2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
